# Konnect Bill Payment Verification Checks added
# - Update the test-case description in A2 to call out the
#   Bene_Verification_pay check explicitly.
# - Column A grows to fit the new (longer) text -> widen/refresh bestfit.
# - Selection cursor left on A9 (where the author's cursor ended up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend the scenario description held in A2.
$ws.Range("A2").Value = "As a user I want to verify Already Added Beneficiaries of Bill Payment Bene_Verification_pay"

# 2. Re-fit column A for the longer text (closest achievable width to the
#    authored 85.140625 "best fit" character width).
$ws.Columns.Item(1).ColumnWidth = 84.33333333333333

# 3. Leave the active selection on A9, matching the saved view state.
$ws.Range("A9").Select()
